$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the MSE / MAE values with the latest computed precision
$ws.Range("B2").Value = 0.0831172165082084
$ws.Range("D2").Value = 0.2072673588334908

# Add the new "Modelo" header in F1, reusing the header formatting (style s="1")
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Add the model description in F2
$ws.Range("F2").Value = "Pipeline(steps=[('model', GradientBoostingRegressor(n_estimators=150))])"
